# Update odds values on the Jogos_da_Semana_FlashScore sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 4.33
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 1.85
$ws.Range("L2").Value = 2.5
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
$ws.Range("AI2").Value = 8.5
$ws.Range("AV2").Value = 51
$ws.Range("BB2").Value = 151

# --- Row 8 updates ---
$ws.Range("G8").Value = 7
$ws.Range("I8").Value = 1.38
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("Z8").Value = 81
$ws.Range("AA8").Value = 51
$ws.Range("AE8").Value = 19
$ws.Range("AJ8").Value = 9
$ws.Range("AK8").Value = 9.5
$ws.Range("AN8").Value = 8.5
$ws.Range("AQ8").Value = 126
$ws.Range("AX8").Value = 6.5
